# Auto-generated update of computed market-price/profit columns (H-N)
# across the 8 Leve sheets, per scheduled price-refresh run.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 614.8
$ws.Range("I2").Value = 418.5
$ws.Range("K2").Value = 418.5
$ws.Range("M2").Value = -305.5
$ws.Range("H62").Value = 19679.883
$ws.Range("I62").Value = 17087.25
$ws.Range("K62").Value = 17087.25
$ws.Range("M62").Value = -16463.25
$ws.Range("H65").Value = 19679.883
$ws.Range("I65").Value = 17087.25
$ws.Range("K65").Value = 85436.25
$ws.Range("M65").Value = -82316.25
$ws.Range("H98").Value = 1070.8667
$ws.Range("J98").Value = 1220
$ws.Range("L98").Value = 1220
$ws.Range("N98").Value = -4216
$ws.Range("H111").Value = 3838.4443
$ws.Range("I111").Value = 2854.8333
$ws.Range("K111").Value = 8564.499899999999
$ws.Range("M111").Value = -5497.499899999999
$ws.Range("H122").Value = 1070.8667
$ws.Range("J122").Value = 1220
$ws.Range("L122").Value = 3660
$ws.Range("N122").Value = -8560

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2978655.5
$ws.Range("I74").Value = 3677509.2
$ws.Range("K74").Value = 3677509.2
$ws.Range("M74").Value = -3676635.2
$ws.Range("H77").Value = 2978655.5
$ws.Range("I77").Value = 3677509.2
$ws.Range("K77").Value = 18387546
$ws.Range("M77").Value = -18383178
$ws.Range("H97").Value = 709.7222
$ws.Range("J97").Value = 128
$ws.Range("L97").Value = 128
$ws.Range("N97").Value = -1120
$ws.Range("H122").Value = 3885.8635
$ws.Range("I122").Value = 3838.3333
$ws.Range("K122").Value = 11514.9999
$ws.Range("M122").Value = -9064.999899999999
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H132").Value = 566697.5
$ws.Range("I132").Value = 660055.7
$ws.Range("J132").Value = 6548.6665
$ws.Range("K132").Value = 1980167.1
$ws.Range("L132").Value = 19645.9995
$ws.Range("M132").Value = -1977637.1
$ws.Range("N132").Value = -24705.9995

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 688294.6
$ws.Range("I134").Value = 1008576
$ws.Range("K134").Value = 3025728
$ws.Range("M134").Value = -3023193
$ws.Range("H141").Value = 70000
$ws.Range("J141").Value = 70000
$ws.Range("L141").Value = 70000
$ws.Range("N141").Value = -80360

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H125").Value = 71663.336
$ws.Range("J125").Value = 71663.336
$ws.Range("L125").Value = 71663.336
$ws.Range("N125").Value = -76583.336
$ws.Range("H132").Value = 27293556
$ws.Range("I132").Value = 29424530
$ws.Range("J132").Value = 21255798
$ws.Range("K132").Value = 88273590
$ws.Range("L132").Value = 63767394
$ws.Range("M132").Value = -88271060
$ws.Range("N132").Value = -63772454
$ws.Range("H134").Value = 3768211
$ws.Range("I134").Value = 22144.875
$ws.Range("J134").Value = 18752476
$ws.Range("K134").Value = 66434.625
$ws.Range("L134").Value = 56257428
$ws.Range("M134").Value = -63899.625
$ws.Range("N134").Value = -56262498

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 250821.75
$ws.Range("I14").Value = 250821.75
$ws.Range("K14").Value = 752465.25
$ws.Range("M14").Value = -752292.25
$ws.Range("H37").Value = 95306.30499999999
$ws.Range("J37").Value = 95306.30499999999
$ws.Range("L37").Value = 285918.915
$ws.Range("N37").Value = -286142.915
$ws.Range("H38").Value = 17.166666
$ws.Range("I38").Value = 14.5
$ws.Range("J38").Value = 22.5
$ws.Range("K38").Value = 43.5
$ws.Range("L38").Value = 67.5
$ws.Range("M38").Value = 303.5
$ws.Range("N38").Value = -761.5
$ws.Range("H107").Value = 595
$ws.Range("J107").Value = 595
$ws.Range("L107").Value = 1785
$ws.Range("N107").Value = -5625

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2278.625
$ws.Range("I80").Value = 2282.15
$ws.Range("J80").Value = 2272.75
$ws.Range("K80").Value = 2282.15
$ws.Range("L80").Value = 2272.75
$ws.Range("M80").Value = -1284.15
$ws.Range("N80").Value = -4268.75
$ws.Range("H83").Value = 2278.625
$ws.Range("I83").Value = 2282.15
$ws.Range("J83").Value = 2272.75
$ws.Range("K83").Value = 11410.75
$ws.Range("L83").Value = 11363.75
$ws.Range("M83").Value = -6418.75
$ws.Range("N83").Value = -21347.75
$ws.Range("H102").Value = 5531.8
$ws.Range("I102").Value = 4914.8335
$ws.Range("J102").Value = 7999.6665
$ws.Range("K102").Value = 4914.8335
$ws.Range("L102").Value = 7999.6665
$ws.Range("M102").Value = -3292.8335
$ws.Range("N102").Value = -11243.6665
$ws.Range("H126").Value = 880099.6
$ws.Range("I126").Value = 1517488.9
$ws.Range("J126").Value = 3689.375
$ws.Range("K126").Value = 4552466.699999999
$ws.Range("L126").Value = 11068.125
$ws.Range("M126").Value = -4549996.699999999
$ws.Range("N126").Value = -16008.125
$ws.Range("H131").Value = 60000
$ws.Range("J131").Value = 60000
$ws.Range("L131").Value = 60000
$ws.Range("N131").Value = -70080
$ws.Range("H132").Value = 225724.4
$ws.Range("I132").Value = 286100.78
$ws.Range("J132").Value = 9375.75
$ws.Range("K132").Value = 858302.3400000001
$ws.Range("L132").Value = 28127.25
$ws.Range("M132").Value = -855772.3400000001
$ws.Range("N132").Value = -33187.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4081.7058
$ws.Range("I7").Value = 3956.3572
$ws.Range("K7").Value = 3956.3572
$ws.Range("M7").Value = -3844.3572
$ws.Range("H40").Value = 4418.5293
$ws.Range("I40").Value = 4316.5386
$ws.Range("J40").Value = 4750
$ws.Range("K40").Value = 4316.5386
$ws.Range("L40").Value = 4750
$ws.Range("M40").Value = -4180.5386
$ws.Range("N40").Value = -5022
$ws.Range("H61").Value = 1886
$ws.Range("I61").Value = 1848
$ws.Range("K61").Value = 1848
$ws.Range("M61").Value = -1646
$ws.Range("H113").Value = 1886
$ws.Range("I113").Value = 1848
$ws.Range("K113").Value = 1848
$ws.Range("M113").Value = 322
$ws.Range("H122").Value = 5091.3076
$ws.Range("I122").Value = 4931.25
$ws.Range("K122").Value = 14793.75
$ws.Range("M122").Value = -12343.75
$ws.Range("H126").Value = 4081.7058
$ws.Range("I126").Value = 3956.3572
$ws.Range("K126").Value = 11869.0716
$ws.Range("M126").Value = -9399.071599999999
$ws.Range("H132").Value = 3871076.8
$ws.Range("J132").Value = 10000
$ws.Range("L132").Value = 30000
$ws.Range("N132").Value = -35060

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 14886.75
$ws.Range("J41").Value = 14886.75
$ws.Range("L41").Value = 14886.75
$ws.Range("N41").Value = -15666.75
$ws.Range("H100").Value = 1806.6
$ws.Range("I100").Value = 1549.6875
$ws.Range("K100").Value = 3099.375
$ws.Range("M100").Value = -2558.375
$ws.Range("H132").Value = 8539399
$ws.Range("I132").Value = 10032359
$ws.Range("K132").Value = 30097077
$ws.Range("M132").Value = -30094547
